# Auto-generated edit script: updates crypto price/volume table to match target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.348.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.388.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.39%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.195'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.75%  '
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.27%  '
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '678.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.927.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.395.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.375.56'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.94%  '
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.31%  '
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '33.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '558.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'OKB'
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '58.63'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.22%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'dogwifhat'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.42%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.669.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.68'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.03%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.139'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.33%  '
$ws.Range("E40").Value = '  +3.72%  '
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("E42").Value = '  +4.56%  '
$ws.Range("E43").Value = '  +1.39%  '
$ws.Range("E44").Value = '  +3.71%  '
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("E48").Value = '  +6.52%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.63'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.06%  '
